# Appends 4 repeated "skip/default" record pairs (rows 172-179) to Sheet1,
# mirroring the existing A:G "profile processed" log pattern already present
# in the sheet (e.g. rows 109/110), but with an extra column H value of
# "Default" on the second row of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$profileUrl = "https://www.instagram.com/_heismannu_"
$profileName = "Koe"
$skip = "Skip"
$default = "Default"

$startRow = 172

for ($i = 0; $i -lt 4; $i++) {
    $row1 = $startRow + ($i * 2)
    $row2 = $row1 + 1

    # First row of the pair: A, B, C, D, E all text
    $ws.Cells.Item($row1, 1).Value = $profileUrl
    $ws.Cells.Item($row1, 2).Value = $profileName
    $ws.Cells.Item($row1, 3).Value = $skip
    $ws.Cells.Item($row1, 4).Value = $skip
    $ws.Cells.Item($row1, 5).Value = $skip

    # Second row of the pair: A, B text; C, F, G boolean FALSE; H text "Default"
    $ws.Cells.Item($row2, 1).Value = $profileUrl
    $ws.Cells.Item($row2, 2).Value = $profileName
    $ws.Cells.Item($row2, 3).Value = $false
    $ws.Cells.Item($row2, 6).Value = $false
    $ws.Cells.Item($row2, 7).Value = $false
    $ws.Cells.Item($row2, 8).Value = $default
}
